$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Insert a new column before the existing "username" column (C) to hold
# browser info, shifting username/password/fname from C/D/E to D/E/F.
$ws.Range("C1").EntireColumn.Insert()

# Header
$ws.Range("C1").Value = "browser"

# Per-row browser values
$ws.Range("C2").Value = "chrome"
$ws.Range("C3").Value = "firefox"
$ws.Range("C4").Value = "chrome"
$ws.Range("C5").Value = "firefox"

# Update the active selection to match the author's edit position
[void]$ws.Range("E3").Select()
